$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "AddCustomerTest"

# Header + data for columns A-C first
$ws.Range("A1").Value = "firstName"
$ws.Range("B1").Value = "lastName"
$ws.Range("C1").Value = "postCode"

$ws.Range("A2").Value = "Anne"
$ws.Range("B2").Value = "Zimmermann"
$ws.Range("C2").Value = "89046-305"

# Column D last
$ws.Range("D1").Value = "alertText"
$ws.Range("D2").Value = "Customer added successfully"

# Style header row: JetBrains Mono 10, vertical center
$headerRange = $ws.Range("A1:D1")
$headerRange.Font.Name = "JetBrains Mono"
$headerRange.Font.Size = 10
$headerRange.Font.Color = 13023145
$headerRange.VerticalAlignment = -4108

# Column widths
$ws.Columns.Item(1).ColumnWidth = 10.6328125
$ws.Columns.Item(2).ColumnWidth = 12
$ws.Columns.Item(3).ColumnWidth = 9.54296875
$ws.Columns.Item(4).ColumnWidth = 39.81640625

# Selection
$ws.Range("E3").Select()

# Window size/position
$excel.ActiveWindow.Left = 2280
$excel.ActiveWindow.Top = 2280
$excel.ActiveWindow.Width = 14400
$excel.ActiveWindow.Height = 7460
